$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1180.2084
$ws.Range("J17").Value = 1180.2084
$ws.Range("L17").Value = 3540.6252
$ws.Range("N17").Value = -3876.6252
$ws.Range("H58").Value = 8710.064
$ws.Range("J58").Value = 10789.2
$ws.Range("L58").Value = 32367.6
$ws.Range("N58").Value = -32667.6
$ws.Range("H64").Value = 6387.4
$ws.Range("I64").Value = 6677.6665
$ws.Range("K64").Value = 6677.6665
$ws.Range("M64").Value = -6429.6665
$ws.Range("H67").Value = 6387.4
$ws.Range("I67").Value = 6677.6665
$ws.Range("K67").Value = 6677.6665
$ws.Range("M67").Value = -5819.6665
$ws.Range("H100").Value = 3958
$ws.Range("J100").Value = 4899.2856
$ws.Range("L100").Value = 4899.2856
$ws.Range("N100").Value = -5981.2856
$ws.Range("H138").Value = 3783.158
$ws.Range("I138").Value = 2869.2354
$ws.Range("J138").Value = 4523
$ws.Range("K138").Value = 8607.706200000001
$ws.Range("L138").Value = 13569
$ws.Range("M138").Value = -3467.706200000001
$ws.Range("N138").Value = -23849
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 1048.4
$ws.Range("I141").Value = 1048.4
$ws.Range("K141").Value = 3145.2
$ws.Range("M141").Value = 2034.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3904.7856
$ws.Range("I61").Value = 2430.5833
$ws.Range("K61").Value = 2430.5833
$ws.Range("M61").Value = -2218.5833
$ws.Range("H63").Value = 7883.486
$ws.Range("J63").Value = 9427
$ws.Range("L63").Value = 9427
$ws.Range("N63").Value = -10799
$ws.Range("H66").Value = 7883.486
$ws.Range("J66").Value = 9427
$ws.Range("L66").Value = 47135
$ws.Range("N66").Value = -53999
$ws.Range("H74").Value = 5825
$ws.Range("I74").Value = 4737.5835
$ws.Range("J74").Value = 7999.8335
$ws.Range("K74").Value = 4737.5835
$ws.Range("L74").Value = 7999.8335
$ws.Range("M74").Value = -3863.5835
$ws.Range("N74").Value = -9747.833500000001
$ws.Range("H77").Value = 5825
$ws.Range("I77").Value = 4737.5835
$ws.Range("J77").Value = 7999.8335
$ws.Range("K77").Value = 23687.9175
$ws.Range("L77").Value = 39999.1675
$ws.Range("M77").Value = -19319.9175
$ws.Range("N77").Value = -48735.1675
$ws.Range("H102").Value = 2142.6216
$ws.Range("I102").Value = 2134.2
$ws.Range("K102").Value = 2134.2
$ws.Range("M102").Value = -512.1999999999998
$ws.Range("H132").Value = 4752.8867
$ws.Range("I132").Value = 2724.55
$ws.Range("K132").Value = 8173.650000000001
$ws.Range("M132").Value = -5643.650000000001
$ws.Range("H136").Value = 3904.7856
$ws.Range("I136").Value = 2430.5833
$ws.Range("K136").Value = 7291.749899999999
$ws.Range("M136").Value = -4741.749899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 8146.778
$ws.Range("I96").Value = 8146.778
$ws.Range("K96").Value = 8146.778
$ws.Range("M96").Value = -5400.778
$ws.Range("H105").Value = 17546584
$ws.Range("J105").Value = 2696.5386
$ws.Range("L105").Value = 2696.5386
$ws.Range("N105").Value = -6190.5386
$ws.Range("H134").Value = 8860.416999999999
$ws.Range("I134").Value = 4950
$ws.Range("J134").Value = 10815.625
$ws.Range("K134").Value = 14850
$ws.Range("L134").Value = 32446.875
$ws.Range("M134").Value = -12315
$ws.Range("N134").Value = -37516.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 55559510
$ws.Range("I31").Value = 125002650
$ws.Range("J31").Value = 5004.7
$ws.Range("K31").Value = 125002650
$ws.Range("L31").Value = 5004.7
$ws.Range("M31").Value = -125002355
$ws.Range("N31").Value = -5594.7
$ws.Range("H34").Value = 55559510
$ws.Range("I34").Value = 125002650
$ws.Range("J34").Value = 5004.7
$ws.Range("K34").Value = 125002650
$ws.Range("L34").Value = 5004.7
$ws.Range("M34").Value = -125002448
$ws.Range("N34").Value = -5408.7
$ws.Range("H62").Value = 12244.682
$ws.Range("I62").Value = 9384
$ws.Range("J62").Value = 13579.667
$ws.Range("K62").Value = 9384
$ws.Range("L62").Value = 13579.667
$ws.Range("M62").Value = -8760
$ws.Range("N62").Value = -14827.667
$ws.Range("H65").Value = 12244.682
$ws.Range("I65").Value = 9384
$ws.Range("J65").Value = 13579.667
$ws.Range("K65").Value = 46920
$ws.Range("L65").Value = 67898.33499999999
$ws.Range("M65").Value = -43800
$ws.Range("N65").Value = -74138.33499999999
$ws.Range("H105").Value = 1888.1666
$ws.Range("I105").Value = 1787.091
$ws.Range("K105").Value = 1787.091
$ws.Range("M105").Value = -40.09099999999989
$ws.Range("H132").Value = 51325.31
$ws.Range("I132").Value = 3725.4092
$ws.Range("J132").Value = 313124.75
$ws.Range("K132").Value = 11176.2276
$ws.Range("L132").Value = 939374.25
$ws.Range("M132").Value = -8646.2276
$ws.Range("N132").Value = -944434.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 2303
$ws.Range("I57").Value = 693.5
$ws.Range("K57").Value = 2080.5
$ws.Range("M57").Value = -1521.5
$ws.Range("H68").Value = 3581.1765
$ws.Range("J68").Value = 4506.9165
$ws.Range("L68").Value = 13520.7495
$ws.Range("N68").Value = -15142.7495
$ws.Range("H71").Value = 3581.1765
$ws.Range("J71").Value = 4506.9165
$ws.Range("L71").Value = 40562.2485
$ws.Range("N71").Value = -48674.2485
$ws.Range("H87").Value = 2124.8333
$ws.Range("I87").Value = 1549.8
$ws.Range("J87").Value = 5000
$ws.Range("K87").Value = 4649.4
$ws.Range("L87").Value = 15000
$ws.Range("M87").Value = -3401.4
$ws.Range("N87").Value = -17496
$ws.Range("H90").Value = 2124.8333
$ws.Range("I90").Value = 1549.8
$ws.Range("J90").Value = 5000
$ws.Range("K90").Value = 13948.2
$ws.Range("L90").Value = 45000
$ws.Range("M90").Value = -7708.199999999999
$ws.Range("N90").Value = -57480
$ws.Range("H97").Value = 538.55554
$ws.Range("I97").Value = 448.6
$ws.Range("J97").Value = 651
$ws.Range("K97").Value = 1345.8
$ws.Range("L97").Value = 1953
$ws.Range("M97").Value = -849.8000000000002
$ws.Range("N97").Value = -2945

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 31255746
$ws.Range("I40").Value = 33338462
$ws.Range("K40").Value = 33338462
$ws.Range("M40").Value = -33338326
$ws.Range("H46").Value = 2329007.5
$ws.Range("I46").Value = 25001500
$ws.Range("J46").Value = 3623.6924
$ws.Range("K46").Value = 25001500
$ws.Range("L46").Value = 3623.6924
$ws.Range("M46").Value = -25001312
$ws.Range("N46").Value = -3999.6924
$ws.Range("H132").Value = 4124.8623
$ws.Range("I132").Value = 2912.5789
$ws.Range("K132").Value = 8737.736699999999
$ws.Range("M132").Value = -6207.736699999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 14074.25
$ws.Range("I96").Value = 1948
$ws.Range("J96").Value = 18116.334
$ws.Range("K96").Value = 1948
$ws.Range("L96").Value = 18116.334
$ws.Range("M96").Value = -575
$ws.Range("N96").Value = -20862.334
$ws.Range("H100").Value = 952.61536
$ws.Range("I100").Value = 665.125
$ws.Range("K100").Value = 1330.25
$ws.Range("M100").Value = -789.25
$ws.Range("H132").Value = 9951.556
$ws.Range("I132").Value = 10078.361
$ws.Range("J132").Value = 9444.333000000001
$ws.Range("K132").Value = 30235.083
$ws.Range("L132").Value = 28332.999
$ws.Range("M132").Value = -27705.083
$ws.Range("N132").Value = -33392.999
$ws.Range("H136").Value = 7561.9165
$ws.Range("I136").Value = 4667.875
$ws.Range("K136").Value = 14003.625
$ws.Range("M136").Value = -11453.625
